# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Row -> new F value
$updates = @{
    2  = 8410
    3  = 7983
    5  = 193
    8  = 135
    9  = 134
    11 = 238
    12 = 719
    13 = 141
    14 = 2035
    16 = 62
    19 = 135
    20 = 38
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
